$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the CodeSystem URL (pythia -> cicada)
$ws.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/EvalStatus"

# Update the generation Date
$ws.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a new "Jurisdiction" row right after "Contact" (row 10), before "Description"
# Copy row 10's formatting down so the new row keeps the same style as its neighbours.
$ws.Rows.Item(11).Insert()
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
